$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "CreateAccountTest" (3rd sheet): selection moved from D5 -> A5
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
[void]$ws3.Range("A5").Select()

# ---------------------------------------------------------------------
# Sheet "LoginTest" (2nd sheet): rework rows 1-5, add hyperlink, etc.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Slightly tighter header / data rows (rows 1-2 keep their old content,
# just pick up a shorter row height).
$ws2.Rows.Item(1).RowHeight = 14.4
$ws2.Rows.Item(2).RowHeight = 14.4

# Row 4 holds what used to live in row 3 (password), carried down with
# its original formatting intact.
$ws2.Range("B3").Copy($ws2.Range("B4")) | Out-Null

# Row 3 becomes a new "forgot password" style row: a hyperlinked e-mail
# address in column A, plain text in B and C.
$ws2.Range("B3").Style = "Normal"
$ws2.Range("C3").Style = "Normal"
$ws2.Range("B3").Value = "test"
$ws2.Range("C3").Value = "Y"
$ws2.Range("A3").Value = "Testug@1test.asu.edu"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:Testug@1test.asu.edu")
$ws2.Range("A3").Style = "Hyperlink"

# Row 4's username cell now carries the same "Hyperlink" visual style used
# elsewhere in the workbook (CreateAccountTest sheet) for the username
# column, plus a trailing Runmode flag.
$ws2.Range("A4").Value = "testug@test.asu.edu"
$ws2.Range("A4").Style = "Hyperlink"
$ws2.Range("C4").Style = "Normal"
$ws2.Range("C4").Value = "Y"

# Row 5: a lone styled, empty cell.
$ws2.Range("A5").Style = "Hyperlink"

$ws2.Rows.Item(3).RowHeight = 15
$ws2.Rows.Item(4).RowHeight = 15
$ws2.Rows.Item(5).RowHeight = 15

[void]$ws2.Range("A10").Select()
